$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Remove the _GoBack bookmark that used to sit after the title
#    ("Matriz de administración del tiempo").
# -----------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# -----------------------------------------------------------------
# 2) Merge the "<Day> <num>" + "/" runs into a single run so the
#    slash is no longer a separate <w:r>.  Doing a plain
#    Find/Replace that crosses the run boundary merges the two
#    runs into one - but it also drops sibling, non-text child
#    elements (e.g. <w:lastRenderedPageBreak/>) that happened to
#    live in the first run when the match spans both runs. So we
#    first delete the lone "/" run's character, then replace the
#    day/number text (a match fully inside the first run) with the
#    day/number text plus the trailing slash - this never crosses a
#    run boundary and therefore keeps any sibling elements intact.
# -----------------------------------------------------------------
function Merge-DaySlash([string]$label) {
    $find = $d.Content
    $find.Find.ClearFormatting()
    $find.Find.Text = "$label/"
    if (-not $find.Find.Execute()) {
        return
    }
    # Delete just the trailing "/" character (the end of the match).
    $slash = $d.Range($find.End - 1, $find.End)
    $slash.Delete()

    # Now replace the (single-run) day/number text with itself plus "/".
    $d.Content.Find.Execute($label, $true, $false, $false, $false, $false, `
                             $true, 1, $false, "$label/", 2) | Out-Null
}

Merge-DaySlash "Jueves 1"
Merge-DaySlash "Viernes 2"

# -----------------------------------------------------------------
# 3) Split "Realizar la tarea de inglés de 8:00 a 9:00am" into three
#    runs, inserting new text and a fresh _GoBack bookmark between
#    the 2nd and 3rd runs.  We rebuild the whole paragraph via
#    InsertXML (Range.InsertXML operates at paragraph granularity
#    here) so we restate the paragraph's own rsid attributes and
#    <w:pPr> verbatim to leave them unchanged.
# -----------------------------------------------------------------
$target = $d.Content
$target.Find.ClearFormatting()
$target.Find.Text = "Realizar la tarea de inglés de 8:00 a 9:00am"
if ($target.Find.Execute()) {
    $sentenceRange = $d.Range($target.Start, $target.End)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes" ?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        '<w:p w:rsidR="00CB15D7" w:rsidRDefault="00CB15D7" w:rsidP="00C77033">' + `
        '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">Realizar la tarea de ingl&#233;s </w:t></w:r>' + `
        '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">y subir el examen de programaci&#243;n </w:t></w:r>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
        '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>de 8:00 a 9:00am</w:t></w:r>' + `
        '</w:p></w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $sentenceRange.InsertXML($xml)
}
